$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 107.5
$ws.Range("I12").Value = 107.5
$ws.Range("K12").Value = 107.5
$ws.Range("M12").Value = 62.5
$ws.Range("H19").Value = 859
$ws.Range("I19").Value = 1090.2
$ws.Range("J19").Value = 570
$ws.Range("K19").Value = 1090.2
$ws.Range("L19").Value = 570
$ws.Range("M19").Value = -915.2
$ws.Range("N19").Value = -920
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H132").Value = 5370.5713
$ws.Range("J132").Value = 8999.5
$ws.Range("L132").Value = 26998.5
$ws.Range("N132").Value = -32058.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3573.4546
$ws.Range("I122").Value = 3379.4
$ws.Range("K122").Value = 10138.2
$ws.Range("M122").Value = -7688.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 87
$ws.Range("I12").Value = 87
$ws.Range("K12").Value = 87
$ws.Range("M12").Value = 81
$ws.Range("H76").Value = 22599.4
$ws.Range("J76").Value = 21499.5
$ws.Range("L76").Value = 21499.5
$ws.Range("N76").Value = -22129.5
$ws.Range("H79").Value = 22599.4
$ws.Range("J79").Value = 21499.5
$ws.Range("L79").Value = 21499.5
$ws.Range("N79").Value = -23683.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 388.66666
$ws.Range("J7").Value = 449.66666
$ws.Range("L7").Value = 449.66666
$ws.Range("N7").Value = -675.66666
$ws.Range("H22").Value = 683.55554
$ws.Range("I22").Value = 854.25
$ws.Range("J22").Value = 547
$ws.Range("K22").Value = 854.25
$ws.Range("L22").Value = 547
$ws.Range("M22").Value = -504.25
$ws.Range("N22").Value = -1247
$ws.Range("H58").Value = 16599.4
$ws.Range("I58").Value = 5999
$ws.Range("J58").Value = 19249.5
$ws.Range("K58").Value = 5999
$ws.Range("L58").Value = 19249.5
$ws.Range("M58").Value = -5796
$ws.Range("N58").Value = -19655.5
$ws.Range("H134").Value = 7841.0713
$ws.Range("I134").Value = 3327.8572
$ws.Range("K134").Value = 9983.5716
$ws.Range("M134").Value = -7448.571599999999
$ws.Range("H136").Value = 16599.4
$ws.Range("I136").Value = 5999
$ws.Range("J136").Value = 19249.5
$ws.Range("K136").Value = 17997
$ws.Range("L136").Value = 57748.5
$ws.Range("M136").Value = -15447
$ws.Range("N136").Value = -62848.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 229.63637
$ws.Range("I2").Value = 313.625
$ws.Range("J2").Value = 5.6666665
$ws.Range("K2").Value = 1881.75
$ws.Range("L2").Value = 33.999999
$ws.Range("M2").Value = -1768.75
$ws.Range("N2").Value = -259.999999
$ws.Range("H4").Value = 49052856
$ws.Range("J4").Value = 73579120
$ws.Range("L4").Value = 220737360
$ws.Range("N4").Value = -220737584
$ws.Range("H36").Value = 2
$ws.Range("I36").Value = 2
$ws.Range("K36").Value = 6
$ws.Range("M36").Value = 163
$ws.Range("H86").Value = 300
$ws.Range("I86").Value = 300
$ws.Range("K86").Value = 900
$ws.Range("M86").Value = 286
$ws.Range("H89").Value = 300
$ws.Range("I89").Value = 300
$ws.Range("K89").Value = 2700
$ws.Range("M89").Value = 3228
$ws.Range("H114").Value = 2193.4
$ws.Range("J114").Value = 2129.625
$ws.Range("L114").Value = 6388.875
$ws.Range("N114").Value = -12896.875
$ws.Range("H117").Value = 3166.6667
$ws.Range("I117").Value = 1000
$ws.Range("K117").Value = 3000
$ws.Range("M117").Value = 442

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 53.692307
$ws.Range("I2").Value = 60.75
$ws.Range("J2").Value = 42.4
$ws.Range("K2").Value = 60.75
$ws.Range("L2").Value = 42.4
$ws.Range("M2").Value = 52.25
$ws.Range("N2").Value = -268.4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H102").Value = 1122.3636
$ws.Range("I102").Value = 1122.3636
$ws.Range("K102").Value = 1122.3636
$ws.Range("M102").Value = 499.6364000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 1150.3334
$ws.Range("I22").Value = 1300
$ws.Range("J22").Value = 851
$ws.Range("K22").Value = 1300
$ws.Range("L22").Value = 851
$ws.Range("M22").Value = -1005
$ws.Range("N22").Value = -1441
$ws.Range("H27").Value = 1150.3334
$ws.Range("I27").Value = 1300
$ws.Range("J27").Value = 851
$ws.Range("K27").Value = 1300
$ws.Range("L27").Value = 851
$ws.Range("M27").Value = -1193
$ws.Range("N27").Value = -1065
$ws.Range("H55").Value = 2000
$ws.Range("I55").Value = 2000
$ws.Range("K55").Value = 2000
$ws.Range("M55").Value = -1827
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 4099.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4099.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12298.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17198.5
$ws.Range("H132").Value = 13612.875
$ws.Range("J132").Value = 19499.75
$ws.Range("L132").Value = 58499.25
$ws.Range("N132").Value = -63559.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 9995
$ws.Range("J11").Value = 9995
$ws.Range("L11").Value = 9995
$ws.Range("N11").Value = -10279
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H97").Value = 27786
$ws.Range("J97").Value = 27786
$ws.Range("L97").Value = 27786
$ws.Range("N97").Value = -29768
$ws.Range("H116").Value = 17500
$ws.Range("J116").Value = 17500
$ws.Range("L116").Value = 17500
$ws.Range("N116").Value = -26678
$ws.Range("H132").Value = 8533.75
$ws.Range("I132").Value = 6711.6665
$ws.Range("K132").Value = 20134.9995
$ws.Range("M132").Value = -17604.9995
